# Apply edits described by the commit diff:
#  - rename several client names (shared strings propagate to both sheets)
#  - rename the employee id
#  - fill in hours/rate/total figures (simulator full-month coverage numbers)

$wb = $excel.ActiveWorkbook
$wsTime = $wb.Worksheets.Item("Weekly Timesheet")
$wsSchema = $wb.Worksheets.Item("Jason Schema")

# --- Client name corrections ---
# These names appear on both the "Weekly Timesheet" sheet (column B) and the
# "Jason Schema" sheet (column D), so update both places explicitly.
$names = @("Prezzano", "Vincent", "Zygmunt", "Ricca", "Varricchio")
for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $i + 2
    $wsTime.Cells.Item($r, 2).Value = $names[$i]
    $wsSchema.Cells.Item($r, 4).Value = $names[$i]
}

# --- Employee ID correction (Jason Schema!B2:B6) ---
for ($r = 2; $r -le 6; $r++) {
    $wsSchema.Cells.Item($r, 2).Value = "emp_emnnysju"
}

# --- Weekly Timesheet: daily Rate (E) / Total (F) for rows 2-6 ---
for ($r = 2; $r -le 6; $r++) {
    $wsTime.Cells.Item($r, 5).Value = 88
    $wsTime.Cells.Item($r, 6).Value = 704
}

# --- Weekly Timesheet: subtotal rows (F8, F11, F13) ---
$wsTime.Range("F8").Value = 3520
$wsTime.Range("F11").Value = 3520
$wsTime.Range("F13").Value = 3520

# --- Jason Schema: per-row Rate (F) / Total (G) for rows 2-6 ---
for ($r = 2; $r -le 6; $r++) {
    $wsSchema.Cells.Item($r, 6).Value = 88
    $wsSchema.Cells.Item($r, 7).Value = 704
}

$wb.Save()
